$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.334.64'
$ws.Range('E2').Value = '  +11.65%  '
$ws.Range('D3').Value = '1.824.80'
$ws.Range('E3').Value = '  +7.73%  '
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.07'
$ws.Range('E5').Value = '  +4.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.545'
$ws.Range('E6').Value = '  +4.26%  '
$ws.Range('E7').Value = '  +0.39%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '31.58'
$ws.Range('E8').Value = '  +2.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.79'
$ws.Range('E9').Value = '  +1.05%  '
$ws.Range('E10').Value = '  +5.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0677'
$ws.Range('E11').Value = '  +8.03%  '
$ws.Range('E12').Value = '  +3.02%  '
$ws.Range('D13').Value = '2.089.25'
$ws.Range('E13').Value = '  +7.81%  '
$ws.Range('D14').Value = '1.842.73'
$ws.Range('E14').Value = '  +8.82%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.647'
$ws.Range('E15').Value = '  +3.86%  '
$ws.Range('D16').Value = '34.359.05'
$ws.Range('E16').Value = '  +11.61%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '10.28'
$ws.Range('E17').Value = '  -3.59%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.36'
$ws.Range('E18').Value = '  +8.46%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.99'
$ws.Range('E19').Value = '  +5.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '259.82'
$ws.Range('E20').Value = '  +4.20%  '
$ws.Range('D21').Value = '0.0₃0753'
$ws.Range('E21').Value = '  +4.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  +0.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.54'
$ws.Range('E23').Value = '  +2.94%  '
$ws.Range('E24').Value = '  +2.26%  '
$ws.Range('E25').Value = '  +0.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '161.37'
$ws.Range('E26').Value = '  +2.46%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.82'
$ws.Range('E27').Value = '  +5.42%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.20'
$ws.Range('E28').Value = '  +6.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.116'
$ws.Range('E29').Value = '  +4.40%  '
$ws.Range('E30').Value = '  +0.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.89'
$ws.Range('E31').Value = '  +11.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0518'
$ws.Range('E32').Value = '  +3.38%  '
$ws.Range('E33').Value = '  +6.87%  '
$ws.Range('E34').Value = '  +8.03%  '
$ws.Range('D35').Value = '1.577.37'
$ws.Range('E35').Value = '  +3.88%  '
$ws.Range('E36').Value = '  +5.82%  '
$ws.Range('E37').Value = '  +4.00%  '
$ws.Range('E38').Value = '  +5.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.632'
$ws.Range('E39').Value = '  +7.86%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '85.04'
$ws.Range('E40').Value = '  +4.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.86'
$ws.Range('E41').Value = '  +5.81%  '
$ws.Range('E42').Value = '  +1.67%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.920'
$ws.Range('E43').Value = '  +7.49%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.14'
$ws.Range('E44').Value = '  +5.52%  '
$ws.Range('E45').Value = '  +3.34%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.07'
$ws.Range('E46').Value = '  +4.58%  '
$ws.Range('D47').Value = '1.981.82'
$ws.Range('E47').Value = '  +8.00%  '
$ws.Range('E48').Value = '  +5.77%  '
$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.00'
$ws.Range('E49').Value = '  +0.28%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '52.90'
$ws.Range('E50').Value = '  +0.51%  '
$ws.Range('D51').Value = '0.0₆0123'
$ws.Range('E51').Value = '  +8.95%  '
